# Update calibration estimates with values re-computed using a standard
# moving-average approach ("Aggiunti file con considerazioni su calcolo
# lambda; aggiornate stime con std moving average").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.004768262808313045
$ws.Range("B3").Value = 0.7688888110379929
$ws.Range("B4").Value = 0.1028686299518762
$ws.Range("B5").Value = "[1.0, 0.023732346766579347, -0.04638569957282715, -0.12020218491831001, -0.12586660649996517, 0.0520734238363353, 0.01673535069328969, 0.06692951596553767, 0.023141806433832124, -0.00747405290850633, -0.06708637295634426, -0.04364949348862001, 0.002529499991237363, 0.02676620867398922, 0.07631326491308171, 0.039307081085521584, -0.0003959883018781687, -0.05228971117715194, -0.0576242126942052, -0.0008429905176588242]"

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1859068900982436
$ws.Range("B5").Value = 0.551211953425346
$ws.Range("B6").Value = 0.06001643437317947
$ws.Range("B7").Value = 0.1757726890812873
$ws.Range("B8").Value = 0.5576344903718747
$ws.Range("B9").Value = 0.06073382125858157
$ws.Range("B10").Value = "[1.0, 0.09675147650289231, 0.009792639305827817, -0.012175979708188032, -0.008085634873128513, -0.005927402853014145, -0.000864025401580957, 0.006108876158650059, 0.0012046032510095667, 0.0019714942549751363, -0.007454028591247468, -0.000958596087635016, -0.006825574044825298, 0.019693520279518453, 0.028233033545550085, 0.00792329356093317, -0.0020193976291577533, -0.005327743454848329, 0.006059428676784363, -0.012123305886651444]"

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = [double]"-6.046801833794415E-05"
$ws.Range("B3").Value = 0.0001291423356933811
$ws.Range("B4").Value = 0.002273644641537301
$ws.Range("B5").Value = 0.9965913068529439
$ws.Range("B6").Value = "[1.0, 0.025508094498297204, -0.036914982248087626, -0.12904891686240694, -0.14105149491945151, 0.0843438872196942, 0.01259839849905599, 0.07936494500292539, 0.022329981128083726, -0.01637738393338926, -0.06412132810500472, -0.050573945908235704, -0.0005161051551699125, 0.030199658979178395, 0.08652411982689583, 0.03639477320946498, -0.010787075196567478, -0.0645899634636009, -0.06832062922466317, 0.008104214349137943]"

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0003445393786604985
$ws.Range("B3").Value = 0.1114476701672416
$ws.Range("B4").Value = 0.04183617997337925
$ws.Range("B5").Value = [double]"6.063967687329516E-09"
$ws.Range("B6").Value = "[1.0, -0.00813828799688224, -0.03312683481797481, -0.12141465038333596, -0.13987891607105132, 0.09056903874282024, 0.008798226647294832, 0.07904513618011906, 0.02139235280722955, -0.01390631423904896, -0.06052381908524808, -0.04883165863270238, -0.0018764185693398835, 0.03167085473665313, 0.08741876965787235, 0.03618079533964133, -0.008050619274894363, -0.06235225105964492, -0.06545145160590261, 0.01272319741547963]"
$ws.Range("B7").Value = -0.003252355184499491

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.005292607249834123
$ws.Range("B3").Value = 0.09928819197218135
$ws.Range("B4").Value = 0.04592924667097358
$ws.Range("B5").Value = [double]"1.235184730071553E-09"
$ws.Range("B6").Value = "[0.9999999999999998, -0.009948737850429597, -0.04575406623915185, -0.11818367795815124, -0.12613559606362357, 0.0569070586715898, 0.01387094502146015, 0.06742980174804601, 0.022988183084669125, -0.005588449984428338, -0.06515209653768768, -0.0436258711427586, 0.0025149658706444364, 0.024025306550922273, 0.07451998200080466, 0.0379273084255292, 0.0003474979810588748, -0.05183119910598594, -0.05810191184165471, -0.0005125876487555238]"
$ws.Range("B7").Value = -0.02445862898434712
$ws.Range("B9").Value = 0.7681058820971721
